$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 00:25"

# Swap country rank: Ecuador overtakes Polonia (rows 34/35)
$ws.Range("A34").Value = "Ecuador"
$ws.Range("A35").Value = "Polonia"

# Swap country rank: Bulgaria overtakes Australia (rows 84/85)
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("A85").Value = "Australia"

# Updated case-count figures (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)

$ws.Range("B4").Value = 8210951
$ws.Range("C4").Value = 60908
$ws.Range("D4").Value = 5310395
$ws.Range("E4").Value = 2677874
$ws.Range("G4").Value = 839
$ws.Range("H4").Value = 222682

$ws.Range("B6").Value = 5169386
$ws.Range("C6").Value = 27888
$ws.Range("D6").Value = 4599446
$ws.Range("E6").Value = 417480
$ws.Range("G6").Value = 681
$ws.Range("H6").Value = 152460

$ws.Range("B9").Value = 949063
$ws.Range("C9").Value = 17096
$ws.Range("D9").Value = 764859
$ws.Range("E9").Value = 158862
$ws.Range("G9").Value = 421
$ws.Range("H9").Value = 25342

$ws.Range("B14").Value = 698184
$ws.Range("C14").Value = 1770
$ws.Range("D14").Value = 628301
$ws.Range("E14").Value = 51574
$ws.Range("G14").Value = 158
$ws.Range("H14").Value = 18309

$ws.Range("D22").Value = 284600
$ws.Range("E22").Value = 54406

$ws.Range("B34").Value = 150360
$ws.Range("C34").Value = 1277
$ws.Range("D34").Value = 128134
$ws.Range("E34").Value = 9920
$ws.Range("G34").Value = 42
$ws.Range("H34").Value = 12306

$ws.Range("B35").Value = 149903
$ws.Range("C35").Value = 8099
$ws.Range("D35").Value = 85588
$ws.Range("E35").Value = 61007
$ws.Range("G35").Value = 91
$ws.Range("H35").Value = 3308

$ws.Range("B46").Value = 105033
$ws.Range("C46").Value = 118
$ws.Range("D46").Value = 98011
$ws.Range("E46").Value = 934
$ws.Range("G46").Value = 11
$ws.Range("H46").Value = 6088

$ws.Range("B57").Value = 76954
$ws.Range("C57").Value = 333
$ws.Range("D57").Value = 73013
$ws.Range("E57").Value = 3652

$ws.Range("B84").Value = 27507
$ws.Range("C84").Value = 914
$ws.Range("D84").Value = 16678
$ws.Range("E84").Value = 9885
$ws.Range("G84").Value = 15
$ws.Range("H84").Value = 944

$ws.Range("B85").Value = 27357
$ws.Range("C85").Value = 16
$ws.Range("D85").Value = 25047
$ws.Range("E85").Value = 1406
$ws.Range("H85").Value = 904

$ws.Range("B108").Value = 10537
$ws.Range("C108").Value = 145
$ws.Range("D108").Value = 8214
$ws.Range("E108").Value = 2250

$ws.Range("B135").Value = 4931
$ws.Range("C135").Value = 48
$ws.Range("D135").Value = 1425
$ws.Range("E135").Value = 3268
$ws.Range("G135").Value = 4
$ws.Range("H135").Value = 238

$ws.Range("B162").Value = 1996
$ws.Range("C162").Value = 13
$ws.Range("D162").Value = 1489
$ws.Range("E162").Value = 456

$ws.Range("B167").Value = 1207
$ws.Range("C167").Value = 2
$ws.Range("E167").Value = 14
